$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows right after existing row 113 (i.e. at position 114),
# which pushes the former rows 114-152 down to become rows 116-154.
$ws.Rows.Item(114).Resize(2).Insert()

# --- New row 114 ---
$ws.Range("A114").Value = 5
$ws.Range("B114").Value = "Macroferia Regional de Talca"
$ws.Range("C114").Value = "Maule"
$ws.Range("D114").Value = 44559
$ws.Range("E114").Value = 7
$ws.Range("F114").Value = 100112024
$ws.Range("G114").Value = "Choclo"
$ws.Range("H114").Value = "Choclero"
$ws.Range("I114").Value = "Primera"
$ws.Range("J114").Value = 20000
$ws.Range("K114").Value = 300
$ws.Range("L114").Value = 300
$ws.Range("M114").Value = 300
$ws.Range("N114").Value = "`$/unidad"
$ws.Range("O114").Value = "Región del Maule"
$ws.Range("P114").Value = 300
$ws.Range("Q114").Value = 1
$ws.Range("R114").Value = "Hortaliza"

# --- New row 115 ---
$ws.Range("A115").Value = 5
$ws.Range("B115").Value = "Macroferia Regional de Talca"
$ws.Range("C115").Value = "Maule"
$ws.Range("D115").Value = 44559
$ws.Range("E115").Value = 7
$ws.Range("F115").Value = 100112024
$ws.Range("G115").Value = "Choclo"
$ws.Range("H115").Value = "Choclero"
$ws.Range("I115").Value = "Segunda"
$ws.Range("J115").Value = 20000
$ws.Range("K115").Value = 200
$ws.Range("L115").Value = 200
$ws.Range("M115").Value = 200
$ws.Range("N115").Value = "`$/unidad"
$ws.Range("O115").Value = "Región del Maule"
$ws.Range("P115").Value = 200
$ws.Range("Q115").Value = 1
$ws.Range("R115").Value = "Hortaliza"
